# Add a new "Partida 8" sheet (multiplayer game data) after the last
# existing "Partida" sheet, fill it with the same T/V/A score layout used
# by the other Partida sheets, make it the active/selected sheet, and
# nudge the number format on one cell (A22) so it carries an explicit
# "apply number format" flag, matching the recorded edit.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

# Insert the new worksheet right after the current last sheet.
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Partida 8"

# Header row: T / V / A
$ws.Range("A1").Value = "T"
$ws.Range("B1").Value = "V"
$ws.Range("C1").Value = "A"

# Score rows 2-29.
$data = @(
  @(15,-5,10),
  @(10,-5,10),
  @(-5,15,15),
  @(15,-5,-10),
  @(10,-5,-5),
  @(30,-5,25),
  @(-5,-10,-15),
  @(-5,-5,15),
  @(-5,25,25),
  @(-5,-5,20),
  @(35,-5,20),
  @(20,-5,-5),
  @(30,35,25),
  @(-5,25,-5),
  @(20,30,-5),
  @(30,45,30),
  @(30,-5,35),
  @(-5,-5,-10),
  @(-5,-5,-5),
  @(15,35,20),
  @(-5,20,-5),
  @(15,-5,20),
  @(-5,15,-5),
  @(15,10,15),
  @(10,-5,-5),
  @(10,15,15),
  @(10,-5,15),
  @(10,-5,10)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Cell A22 picks up an explicit (re-applied) "General" number format.
$ws.Range("A22").NumberFormat = "General"

# Row 18 carries a slightly taller, explicit row height.
$ws.Rows.Item(18).RowHeight = 16.5

# Put the selection where it ended up in the recorded session and make
# this new sheet the active/selected tab.
$ws.Range("F28").Select()
$ws.Activate()
